$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512554836673157"
$ws1.Range("B2").Value = "go_stims-16512554836369507.csv"
$ws1.Range("B3").Value = "GNG_stims-16512554836505563.csv"
$ws1.Range("B4").Value = "go_stims-1651255483651596.csv"
$ws1.Range("B5").Value = "GNG_stims-16512554836663523.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16512554860199082"
$ws2.Range("B2").Value = "OB-16512554849356887.csv"
$ws2.Range("B3").Value = "TB-16512554853382897.csv"
$ws2.Range("B4").Value = "ZB-match_6-1651255483817136.csv"
$ws2.Range("B5").Value = "ZB-match_0-16512554838441272.csv"
$ws2.Range("B6").Value = "TB-16512554851928227.csv"
$ws2.Range("B7").Value = "ZB-match_2-16512554844141412.csv"
$ws2.Range("B8").Value = "OB-16512554848904274.csv"
$ws2.Range("B9").Value = "OB-16512554847173853.csv"
$ws2.Range("B10").Value = "TB-1651255485995942.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512554860209057"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16512554860839047"
$ws4.Range("B2").Value = "MM_stims-16512554860359385.csv"
$ws4.Range("B3").Value = "ZM_stims-16512554860229044.csv"
$ws4.Range("B4").Value = "MM_stims-16512554860669394.csv"
$ws4.Range("B5").Value = "ZM_stims-16512554860369048.csv"
$ws4.Range("B6").Value = "MM_stims-1651255486082906.csv"
$ws4.Range("B7").Value = "ZM_stims-16512554860679066.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512554861469405"
$ws5.Range("B2").Value = "vSAT_stims-16512554861309469.csv"
$ws5.Range("B3").Value = "SAT_stims-16512554860989392.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512554861159108.csv"
$ws5.Range("B5").Value = "SAT_stims-16512554860869076.csv"
